# "fixed all list and new tag associations"
#
# The List6/List7/List8 "list" paragraph styles were mistakenly based on
# Heading4 (pulling in bold + an outline level used for TOC purposes).
# This cleans that up: drop the Heading4 basis and the contextual spacing
# that came along with the old (wrong) lineage, and restate the
# run-formatting explicitly on the styles themselves (Times New Roman,
# 12pt, with List7/List7Char kept bold like the surrounding "change"
# styles). List1 / List3change / List4change also drop the now-redundant
# contextual-spacing flag.

$d = $word.ActiveDocument

# --- List1: remove contextual spacing -------------------------------------
$d.Styles("List1").NoSpaceBetweenParagraphsOfSameStyle = $false

# --- List6: no longer based on Heading4 ------------------------------------
$s6 = $d.Styles("List6")
$s6.BaseStyle = ""
$s6.NoSpaceBetweenParagraphsOfSameStyle = $false
$s6.Font.Name = "Times New Roman"
$s6.Font.NameAscii = "Times New Roman"
$s6.Font.NameBi = "Times New Roman"
$s6.Font.NameOther = "Times New Roman"
$s6.Font.Size = 12
$s6.Font.SizeBi = 12

# --- List7: no longer based on Heading4 (stays bold) -----------------------
$s7 = $d.Styles("List7")
$s7.BaseStyle = ""
$s7.Font.Name = "Times New Roman"
$s7.Font.NameAscii = "Times New Roman"
$s7.Font.NameBi = "Times New Roman"
$s7.Font.NameOther = "Times New Roman"
$s7.Font.Size = 12
$s7.Font.SizeBi = 12
$s7.Font.Bold = $true

$s7c = $d.Styles("List7Char")
$s7c.Font.Bold = $true

# --- List8: no longer based on Heading4 ------------------------------------
$s8 = $d.Styles("List8")
$s8.BaseStyle = ""
$s8.NoSpaceBetweenParagraphsOfSameStyle = $false
$s8.Font.Name = "Times New Roman"
$s8.Font.NameAscii = "Times New Roman"
$s8.Font.NameBi = "Times New Roman"
$s8.Font.NameOther = "Times New Roman"
$s8.Font.Size = 12
$s8.Font.SizeBi = 12

# --- List3change / List4change: remove contextual spacing ------------------
$d.Styles("List3change").NoSpaceBetweenParagraphsOfSameStyle = $false
$d.Styles("List4change").NoSpaceBetweenParagraphsOfSameStyle = $false

Write-Output "style fixups applied"
